$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date in column C was refreshed for every
# existing data row (2-330): 45204 -> 45205 (2023-10-05 -> 2023-10-06).
$ws.Range("C2:C330").Value = 45205

# Row 330 picks up an explicit row height in the new export.
$ws.Rows.Item(330).RowHeight = 15

# A brand-new cutting notification was appended as row 331.
$ws.Cells.Item(331, 1).Value = "A 47851-2023"
$ws.Cells.Item(331, 2).Value = 45204
$ws.Cells.Item(331, 3).Value = 45205
$ws.Cells.Item(331, 4).Value = "BLEKINGE LÄN"
$ws.Cells.Item(331, 5).Value = "KARLSHAMN"
$ws.Cells.Item(331, 7).Value = 6.8
$ws.Cells.Item(331, 8).Value = 0
$ws.Cells.Item(331, 9).Value = 0
$ws.Cells.Item(331, 10).Value = 0
$ws.Cells.Item(331, 11).Value = 0
$ws.Cells.Item(331, 12).Value = 0
$ws.Cells.Item(331, 13).Value = 0
$ws.Cells.Item(331, 14).Value = 0
$ws.Cells.Item(331, 15).Value = 0
$ws.Cells.Item(331, 16).Value = 0
$ws.Cells.Item(331, 17).Value = 0

# Match the date-formatted style used by column B/C elsewhere, and the
# wrap-text style used by column R elsewhere.
$ws.Range("B331:C331").NumberFormat = "YYYY-MM-DD"
$ws.Range("R331").WrapText = $true
